$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Cells.Item(28, 8).Value = 572.5
$ws.Cells.Item(28, 9).Value = 618.5714
$ws.Cells.Item(28, 10).Value = 250
$ws.Cells.Item(28, 11).Value = 618.5714
$ws.Cells.Item(28, 12).Value = 250
$ws.Cells.Item(28, 13).Value = -133.5714
$ws.Cells.Item(28, 14).Value = -1220
# Row 43 (Leve Item ID 5472)
$ws.Cells.Item(43, 8).Value = 900
$ws.Cells.Item(43, 9).Value = 900
$ws.Cells.Item(43, 11).Value = 900
$ws.Cells.Item(43, 13).Value = -831
# Row 70 (Leve Item ID 12604)
$ws.Cells.Item(70, 8).Value = 51169.2
$ws.Cells.Item(70, 10).Value = 1181.4117
$ws.Cells.Item(70, 12).Value = 3544.2351
$ws.Cells.Item(70, 14).Value = -4084.2351
# Row 73 (Leve Item ID 12604)
$ws.Cells.Item(73, 8).Value = 51169.2
$ws.Cells.Item(73, 10).Value = 1181.4117
$ws.Cells.Item(73, 12).Value = 3544.2351
$ws.Cells.Item(73, 14).Value = -5416.2351
# Row 98 (Leve Item ID 36237)
$ws.Cells.Item(98, 8).Value = 5354.9614
$ws.Cells.Item(98, 9).Value = 3401.3809
$ws.Cells.Item(98, 10).Value = 13560
$ws.Cells.Item(98, 11).Value = 3401.3809
$ws.Cells.Item(98, 12).Value = 13560
$ws.Cells.Item(98, 13).Value = -1903.3809
$ws.Cells.Item(98, 14).Value = -16556
# Row 107 (Leve Item ID 27766)
$ws.Cells.Item(107, 8).Value = 367.53125
$ws.Cells.Item(107, 9).Value = 384.42307
$ws.Cells.Item(107, 10).Value = 294.33334
$ws.Cells.Item(107, 11).Value = 384.42307
$ws.Cells.Item(107, 12).Value = 294.33334
$ws.Cells.Item(107, 13).Value = 1535.57693
$ws.Cells.Item(107, 14).Value = -4134.33334
# Row 113 (Leve Item ID 27775)
$ws.Cells.Item(113, 8).Value = 2750.625
$ws.Cells.Item(113, 9).Value = 2351.25
$ws.Cells.Item(113, 10).Value = 3150
$ws.Cells.Item(113, 11).Value = 2351.25
$ws.Cells.Item(113, 12).Value = 3150
$ws.Cells.Item(113, 13).Value = 902.75
$ws.Cells.Item(113, 14).Value = -9658
# Row 122 (Leve Item ID 36237)
$ws.Cells.Item(122, 8).Value = 5354.9614
$ws.Cells.Item(122, 9).Value = 3401.3809
$ws.Cells.Item(122, 10).Value = 13560
$ws.Cells.Item(122, 11).Value = 10204.1427
$ws.Cells.Item(122, 12).Value = 40680
$ws.Cells.Item(122, 13).Value = -7754.1427
$ws.Cells.Item(122, 14).Value = -45580
# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, 8).Value = 2533.75
$ws.Cells.Item(132, 9).Value = 2088.5833
$ws.Cells.Item(132, 11).Value = 6265.749899999999
$ws.Cells.Item(132, 13).Value = -3735.749899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Cells.Item(2, 8).Value = 717.2549
$ws.Cells.Item(2, 9).Value = 471.6154
$ws.Cells.Item(2, 10).Value = 1515.5834
$ws.Cells.Item(2, 11).Value = 471.6154
$ws.Cells.Item(2, 12).Value = 1515.5834
$ws.Cells.Item(2, 13).Value = -358.6154
$ws.Cells.Item(2, 14).Value = -1741.5834
# Row 64 (Leve Item ID 10664)
$ws.Cells.Item(64, 8).Value = 30000
$ws.Cells.Item(64, 10).Value = 30000
$ws.Cells.Item(64, 12).Value = 30000
$ws.Cells.Item(64, 14).Value = -30496
# Row 67 (Leve Item ID 10664)
$ws.Cells.Item(67, 8).Value = 30000
$ws.Cells.Item(67, 10).Value = 30000
$ws.Cells.Item(67, 12).Value = 30000
$ws.Cells.Item(67, 14).Value = -31716
# Row 80 (Leve Item ID 10667)
$ws.Cells.Item(80, 8).Value = 35933.332
$ws.Cells.Item(80, 10).Value = 35933.332
$ws.Cells.Item(80, 12).Value = 35933.332
$ws.Cells.Item(80, 14).Value = -37929.332
# Row 83 (Leve Item ID 10667)
$ws.Cells.Item(83, 8).Value = 35933.332
$ws.Cells.Item(83, 10).Value = 35933.332
$ws.Cells.Item(83, 12).Value = 107799.996
$ws.Cells.Item(83, 14).Value = -117783.996
# Row 104 (Leve Item ID 18672)
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).ClearContents()
# Row 107 (Leve Item ID 25645)
$ws.Cells.Item(107, 8).Value = 23750
$ws.Cells.Item(107, 10).Value = 23750
$ws.Cells.Item(107, 12).Value = 23750
$ws.Cells.Item(107, 14).Value = -31430
# Row 110 (Leve Item ID 27708)
$ws.Cells.Item(110, 8).Value = 1425.9166
$ws.Cells.Item(110, 9).Value = 1306.7778
$ws.Cells.Item(110, 11).Value = 1306.7778
$ws.Cells.Item(110, 13).Value = 738.2221999999999
# Row 116 (Leve Item ID 27713)
$ws.Cells.Item(116, 8).Value = 717.2549
$ws.Cells.Item(116, 9).Value = 471.6154
$ws.Cells.Item(116, 10).Value = 1515.5834
$ws.Cells.Item(116, 11).Value = 471.6154
$ws.Cells.Item(116, 12).Value = 1515.5834
$ws.Cells.Item(116, 13).Value = 1822.3846
$ws.Cells.Item(116, 14).Value = -6103.5834

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Cells.Item(3, 8).Value = 717.2549
$ws.Cells.Item(3, 9).Value = 471.6154
$ws.Cells.Item(3, 10).Value = 1515.5834
$ws.Cells.Item(3, 11).Value = 471.6154
$ws.Cells.Item(3, 12).Value = 1515.5834
$ws.Cells.Item(3, 13).Value = -357.6154
$ws.Cells.Item(3, 14).Value = -1743.5834
# Row 107 (Leve Item ID 27706)
$ws.Cells.Item(107, 8).Value = 66825
$ws.Cells.Item(107, 9).Value = 66825
$ws.Cells.Item(107, 11).Value = 66825
$ws.Cells.Item(107, 13).Value = -64905
# Row 132 (Leve Item ID 41855)
$ws.Cells.Item(132, 8).Value = 72374.14
$ws.Cells.Item(132, 10).Value = 72374.14
$ws.Cells.Item(132, 12).Value = 72374.14
$ws.Cells.Item(132, 14).Value = -82494.14

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (Leve Item ID 44019)
$ws.Cells.Item(132, 8).Value = 2211.7932
$ws.Cells.Item(132, 9).Value = 1805.76
$ws.Cells.Item(132, 10).Value = 4749.5
$ws.Cells.Item(132, 11).Value = 5417.28
$ws.Cells.Item(132, 12).Value = 14248.5
$ws.Cells.Item(132, 13).Value = -2887.28
$ws.Cells.Item(132, 14).Value = -19308.5

$ws = $wb.Worksheets.Item("CUL")
# Row 127 (Leve Item ID 38263)
$ws.Cells.Item(127, 8).Value = 450
$ws.Cells.Item(127, 10).Value = 450
$ws.Cells.Item(127, 12).Value = 1350
$ws.Cells.Item(127, 14).Value = -11270
# Row 130 (Leve Item ID 36058)
$ws.Cells.Item(130, 8).Value = 2312.5
$ws.Cells.Item(130, 10).Value = 2750
$ws.Cells.Item(130, 12).Value = 8250
$ws.Cells.Item(130, 14).Value = -18290
# Row 131 (Leve Item ID 36060)
$ws.Cells.Item(131, 8).Value = 852.04
$ws.Cells.Item(131, 10).Value = 886.5543
$ws.Cells.Item(131, 12).Value = 2659.6629
$ws.Cells.Item(131, 14).Value = -12739.6629

$ws = $wb.Worksheets.Item("GSM")
# Row 57 (Leve Item ID 2876)
$ws.Cells.Item(57, 8).Value = 15005
$ws.Cells.Item(57, 9).Value = 15005
$ws.Cells.Item(57, 11).Value = 15005
$ws.Cells.Item(57, 13).Value = -14185
# Row 63 (Leve Item ID 11048)
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).ClearContents()
# Row 66 (Leve Item ID 11048)
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).ClearContents()
# Row 97 (Leve Item ID 19940)
$ws.Cells.Item(97, 8).Value = 29324.445
$ws.Cells.Item(97, 9).Value = 29324.445
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 29324.445
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -28828.445
$ws.Cells.Item(97, 14).ClearContents()
# Row 113 (Leve Item ID 27710)
$ws.Cells.Item(113, 8).Value = 992.8461
$ws.Cells.Item(113, 9).Value = 916.1
$ws.Cells.Item(113, 10).Value = 1248.6666
$ws.Cells.Item(113, 11).Value = 916.1
$ws.Cells.Item(113, 12).Value = 1248.6666
$ws.Cells.Item(113, 13).Value = 1253.9
$ws.Cells.Item(113, 14).Value = -5588.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 55 (Leve Item ID 5284)
$ws.Cells.Item(55, 8).Value = 289.5
$ws.Cells.Item(55, 10).Value = 225
$ws.Cells.Item(55, 12).Value = 225
$ws.Cells.Item(55, 14).Value = -571
# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 6463.9546
$ws.Cells.Item(132, 9).Value = 6616.2104
$ws.Cells.Item(132, 10).Value = 5499.6665
$ws.Cells.Item(132, 11).Value = 19848.6312
$ws.Cells.Item(132, 12).Value = 16498.9995
$ws.Cells.Item(132, 13).Value = -17318.6312
$ws.Cells.Item(132, 14).Value = -21558.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 16 (Leve Item ID 26304)
$ws.Cells.Item(16, 8).Value = 35568
$ws.Cells.Item(16, 10).Value = 35568
$ws.Cells.Item(16, 12).Value = 35568
$ws.Cells.Item(16, 14).Value = -36152
# Row 64 (Leve Item ID 11036)
$ws.Cells.Item(64, 8).Value = 23114
$ws.Cells.Item(64, 10).Value = 23114
$ws.Cells.Item(64, 12).Value = 23114
$ws.Cells.Item(64, 14).Value = -23610
# Row 67 (Leve Item ID 11036)
$ws.Cells.Item(67, 8).Value = 23114
$ws.Cells.Item(67, 10).Value = 23114
$ws.Cells.Item(67, 12).Value = 23114
$ws.Cells.Item(67, 14).Value = -24830

Write-Host "Applied Asura_Profits market price updates"
